$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for the new columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the style of an existing header cell (e.g. AC1) to the new header cells
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in Wins/Losses/Ties values for all data rows (2-57)
for ($r = 2; $r -le 57; $r++) {
    $ws.Cells.Item($r, 30).Value = 90   # AD
    $ws.Cells.Item($r, 31).Value = 72   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
